# Daily attendance processing - 2026-01-02 17:57:24
# Swap the order of the "Recorded By" names from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
